$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New invoice entry in row 7: STM32F051 dev board PCB, bought from JLCPCB,
# paid for/received by JM.

# Date (2020-11-23), formatted like the other date cells in column B
# (built-in short-date format, centered horizontally and vertically).
$ws.Range("B7").Value = 44158
$ws.Range("B7").NumberFormat = "mm-dd-yy"
$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("B7").VerticalAlignment = -4108

$ws.Range("C7").Value = "STM32F051 dev board PCB"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 16.26
$ws.Range("F7").Value = "JLCPCB"
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = "JM"

# Match the last-saved selection recorded in the workbook.
$ws.Range("G6").Select()
